$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $addr, $val) {
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws 'D2' '29.229.68'
Set-TextCell $ws 'E2' '  -0.21%  '

# Row 3
Set-TextCell $ws 'D3' '1.840.21'
Set-TextCell $ws 'E3' '  +0.02%  '

# Row 4
Set-TextCell $ws 'D4' '0.9997'
Set-TextCell $ws 'E4' '  +0.02%  '

# Row 5
Set-TextCell $ws 'D5' '240.76'
Set-TextCell $ws 'E5' '  -0.95%  '

# Row 6
Set-TextCell $ws 'D6' '0.6691'
Set-TextCell $ws 'E6' '  -2.36%  '

# Row 7
Set-TextCell $ws 'D7' '1.001'
Set-TextCell $ws 'E7' '  +0.02%  '

# Row 8
Set-TextCell $ws 'D8' '0.07417'
Set-TextCell $ws 'E8' '  -1.06%  '

# Row 9
Set-TextCell $ws 'D9' '0.2960'
Set-TextCell $ws 'E9' '  -2.20%  '

# Row 10
Set-TextCell $ws 'D10' '22.83'
Set-TextCell $ws 'E10' '  -1.52%  '

# Row 11
Set-TextCell $ws 'D11' '0.07711'
Set-TextCell $ws 'E11' '  +0.80%  '

# Row 12
Set-TextCell $ws 'D12' '5.027'
Set-TextCell $ws 'E12' '  -0.86%  '

# Row 13
Set-TextCell $ws 'D13' '0.6787'
Set-TextCell $ws 'E13' '  -0.82%  '

# Row 14
Set-TextCell $ws 'D14' '1.749.23'
Set-TextCell $ws 'E14' '  -4.91%  '

# Row 15
Set-TextCell $ws 'D15' '86.21'
Set-TextCell $ws 'E15' '  -3.32%  '

# Row 16
Set-TextCell $ws 'D16' '6.201'
Set-TextCell $ws 'E16' '  -1.32%  '

# Row 17
Set-TextCell $ws 'B17' 'ShibaInu'
Set-TextCell $ws 'C17' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell $ws 'D17' '0.000008251'
Set-TextCell $ws 'E17' '  +0.22%  '

# Row 18
Set-TextCell $ws 'B18' 'WrappedBTC'
Set-TextCell $ws 'C18' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell $ws 'D18' '28.762.92'
Set-TextCell $ws 'E18' '  -1.69%  '

# Row 19
Set-TextCell $ws 'D19' '228.85'
Set-TextCell $ws 'E19' '  -2.25%  '

# Row 20
Set-TextCell $ws 'D20' '12.54'
Set-TextCell $ws 'E20' '  -0.17%  '

# Row 21
Set-TextCell $ws 'D21' '0.9998'
Set-TextCell $ws 'E21' '  -0.03%  '

# Row 22
Set-TextCell $ws 'D22' '7.231'
Set-TextCell $ws 'E22' '  -2.98%  '

# Row 23
Set-TextCell $ws 'D23' '1.001'
Set-TextCell $ws 'E23' '  +0.06%  '

# Row 24
Set-TextCell $ws 'D24' '160.17'
Set-TextCell $ws 'E24' '  +0.29%  '

# Row 25
Set-TextCell $ws 'D25' '8.713'
Set-TextCell $ws 'E25' '  -1.23%  '

# Row 26
Set-TextCell $ws 'D26' '0.1415'
Set-TextCell $ws 'E26' '  -2.78%  '

# Row 27
Set-TextCell $ws 'D27' '18.04'
Set-TextCell $ws 'E27' '  -0.18%  '

# Row 28
Set-TextCell $ws 'D28' '1.506'
Set-TextCell $ws 'E28' '  -0.82%  '

# Row 29
Set-TextCell $ws 'D29' '4.205'
Set-TextCell $ws 'E29' '  -0.24%  '

# Row 30
Set-TextCell $ws 'D30' '4.081'
Set-TextCell $ws 'E30' '  -1.08%  '

# Row 31
Set-TextCell $ws 'D31' '1.196'
Set-TextCell $ws 'E31' '  -0.46%  '

# Row 32
Set-TextCell $ws 'D32' '0.05359'
Set-TextCell $ws 'E32' '  +4.52%  '

# Row 33
Set-TextCell $ws 'D33' '0.7581'
Set-TextCell $ws 'E33' '  -1.33%  '

# Row 34
Set-TextCell $ws 'D34' '1.870'
Set-TextCell $ws 'E34' '  +1.55%  '

# Row 35
Set-TextCell $ws 'D35' '1.136'
Set-TextCell $ws 'E35' '  +0.04%  '

# Row 36
Set-TextCell $ws 'E36' '  +0.41%  '

# Row 37
Set-TextCell $ws 'D37' '1.333.09'
Set-TextCell $ws 'E37' '  +3.56%  '

# Row 38
Set-TextCell $ws 'D38' '0.01802'
Set-TextCell $ws 'E38' '  -2.20%  '

# Row 39
Set-TextCell $ws 'D39' '2.730'
Set-TextCell $ws 'E39' '  +1.21%  '

# Row 40
Set-TextCell $ws 'D40' '0.9217'
Set-TextCell $ws 'E40' '  -1.96%  '

# Row 41
Set-TextCell $ws 'D41' '5.997'
Set-TextCell $ws 'E41' '  +6.45%  '

# Row 42
Set-TextCell $ws 'D42' '1.0000'
Set-TextCell $ws 'E42' '  -0.02%  '

# Row 43
Set-TextCell $ws 'D43' '103.23'
Set-TextCell $ws 'E43' '  -2.01%  '

# Row 44
Set-TextCell $ws 'B44' 'BabyDogeCoin'
Set-TextCell $ws 'C44' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell $ws 'D44' '0.00000000125'
Set-TextCell $ws 'E44' '  +2.28%  '

# Row 45
Set-TextCell $ws 'B45' 'XinFinNetwork'
Set-TextCell $ws 'C45' 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
Set-TextCell $ws 'D45' '0.07884'
Set-TextCell $ws 'E45' '  +13.68%  '

# Row 46
Set-TextCell $ws 'D46' '0.5164'
Set-TextCell $ws 'E46' '  -0.77%  '

# Row 47
Set-TextCell $ws 'B47' 'Aave'
Set-TextCell $ws 'C47' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell $ws 'D47' '63.87'
Set-TextCell $ws 'E47' '  +1.62%  '

# Row 48
Set-TextCell $ws 'B48' 'RenderToken'
Set-TextCell $ws 'C48' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell $ws 'D48' '1.765'
Set-TextCell $ws 'E48' '  +0.39%  '

# Row 49
Set-TextCell $ws 'B49' 'EnergySwap'
Set-TextCell $ws 'C49' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell $ws 'D49' '9.287'
Set-TextCell $ws 'E49' '  -3.71%  '

# Row 50
Set-TextCell $ws 'B50' 'RocketPoolETH'
Set-TextCell $ws 'C50' 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextCell $ws 'D50' '1.883.69'
Set-TextCell $ws 'E50' '  -5.38%  '

# Row 51
Set-TextCell $ws 'D51' '0.05937'
Set-TextCell $ws 'E51' '  +0.25%  '
